$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271, shifting existing rows 271:343 down to 272:344
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new record's data
$ws.Cells.Item(271, 1).Value = 7
$ws.Cells.Item(271, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(271, 3).Value = 'Ñuble'
$ws.Cells.Item(271, 4).Value = 44551
$ws.Cells.Item(271, 5).Value = 16
$ws.Cells.Item(271, 6).Value = 100112020
$ws.Cells.Item(271, 7).Value = 'Tomate'
$ws.Cells.Item(271, 8).Value = 'Larga vida'
$ws.Cells.Item(271, 9).Value = 'Primera'
$ws.Cells.Item(271, 10).Value = 2400
$ws.Cells.Item(271, 11).Value = 5500
$ws.Cells.Item(271, 12).Value = 6000
$ws.Cells.Item(271, 13).Value = 5750
$ws.Cells.Item(271, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(271, 15).Value = 'Región del Maule'
$ws.Cells.Item(271, 16).Value = 383
$ws.Cells.Item(271, 17).Value = 15
$ws.Cells.Item(271, 18).Value = 'Hortaliza'
